# Amazon_Cucumber - Addaddress_Amazon.xlsx update
# Change the City value in row 2 (D2) from "LA" to "COMMERCE"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AmazonAddress")

$ws.Range("D2").Value = "COMMERCE"

# Reflect the final cell selection recorded in the saved file (K22)
$ws.Range("K22").Select()
